$wb = $excel.ActiveWorkbook

# --- Sheet "DataCombined" (1st sheet) ---------------------------------
$wsData = $wb.Worksheets.Item(1)

# Write D4 (PopulationScenario) before A4/G4 (AciclovirPop) so that the
# shared-string table gets "PopulationScenario" at index 58 and
# "AciclovirPop" at index 59, matching the target workbook layout.
$wsData.Range("D4").Value = "PopulationScenario"

# Row 4: new simulated entry for AciclovirPop
$wsData.Range("A4").Value = "AciclovirPop"
$wsData.Range("B4").Value = "simulated"
$wsData.Range("C4").Value = "Aciclovir simulated"
$wsData.Range("E4").Value = "Organism|PeripheralVenousBlood|Aciclovir|Plasma (Peripheral Venous Blood)"
$wsData.Range("G4").Value = "AciclovirPop"
$wsData.Range("H4").Value = 1
$wsData.Range("I4").Value = "h"

# Row 5: new observed entry for AciclovirPop
$wsData.Range("A5").Value = "AciclovirPop"
$wsData.Range("B5").Value = "observed"
# New row 5 would otherwise inherit column B's "vertical-center" style
# (used only for the header row); reset it back to Normal to match row 2-4.
$wsData.Range("B5").Style = "Normal"
$wsData.Range("C5").Value = "Aciclovri observed"
$wsData.Range("F5").Value = "Laskin 1982.Group A_Aciclovir_1_Human_PeripheralVenousBlood_Plasma_2.5 mg/kg_iv_"
$wsData.Range("G5").Value = "AciclovirPop"
$wsData.Range("H5").Value = 1
$wsData.Range("I5").Value = "min"

# Resize column F to its new (slightly wider, no-longer-autofit) width.
$wsData.Columns.Item(6).ColumnWidth = 80.45182291666667

# --- Sheet "plotConfiguration" (2nd sheet) -----------------------------
$wsPlotConfig = $wb.Worksheets.Item(2)

# Row 5 (plot P4) now references the new population data combination and
# plots the "population" plot type instead of "individual".
$wsPlotConfig.Range("B5").Value = "AciclovirPop"
$wsPlotConfig.Range("C5").Value = "population"

# --- Active sheet / selection bookkeeping ------------------------------
# In the edited workbook, "DataCombined" becomes the active tab with D4
# selected, while "plotConfiguration" is no longer the active tab and has
# B5 selected.
$wsPlotConfig.Activate()
$wsPlotConfig.Range("B5").Select()

$wsData.Activate()
$wsData.Range("D4").Select()
